$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 198, pushing the existing rows 198:326 down to 199:327.
$ws.Rows("198:198").Insert()

# Populate the newly inserted row 198 with its data.
$ws.Cells.Item(198, 1).Value = 4
$ws.Cells.Item(198, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(198, 3).Value = "Los Lagos"
$ws.Cells.Item(198, 4).Value = Get-Date -Year 2022 -Month 7 -Day 25 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item(198, 5).Value = 10
$ws.Cells.Item(198, 6).Value = 100114014
$ws.Cells.Item(198, 7).Value = "Betarraga"
$ws.Cells.Item(198, 8).Value = "Sin especificar"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 500
$ws.Cells.Item(198, 11).Value = 1200
$ws.Cells.Item(198, 12).Value = 1200
$ws.Cells.Item(198, 13).Value = 1200
$ws.Cells.Item(198, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(198, 15).Value = "Región del Maule"
$ws.Cells.Item(198, 16).Value = 240
$ws.Cells.Item(198, 17).Value = 5
$ws.Cells.Item(198, 18).Value = "Hortaliza"
